$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 724.75
$ws.Range("I19").Value = 850
$ws.Range("J19").Value = 599.5
$ws.Range("K19").Value = 850
$ws.Range("L19").Value = 599.5
$ws.Range("M19").Value = -675
$ws.Range("N19").Value = -949.5
$ws.Range("H95").Value = 37166
$ws.Range("J95").Value = 37166
$ws.Range("L95").Value = 37166
$ws.Range("N95").Value = -42658
$ws.Range("H124").Value = 46776
$ws.Range("J124").Value = 46776
$ws.Range("L124").Value = 46776
$ws.Range("N124").Value = -56596
$ws.Range("H128").Value = 48832.75
$ws.Range("J128").Value = 48832.75
$ws.Range("L128").Value = 48832.75
$ws.Range("N128").Value = -58792.75
$ws.Range("H130").Value = 49776
$ws.Range("J130").Value = 49776
$ws.Range("L130").Value = 49776
$ws.Range("N130").Value = -59816
$ws.Range("H137").Value = 3749.1667
$ws.Range("I137").Value = 977.6667
$ws.Range("K137").Value = 2933.0001
$ws.Range("M137").Value = -383.0001000000002

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2205.4119
$ws.Range("I2").Value = 2466.5715
$ws.Range("J2").Value = 986.6667
$ws.Range("K2").Value = 2466.5715
$ws.Range("L2").Value = 986.6667
$ws.Range("M2").Value = -2353.5715
$ws.Range("N2").Value = -1212.6667
$ws.Range("H45").Value = 1809.05
$ws.Range("I45").Value = 1776.9231
$ws.Range("J45").Value = 1868.7142
$ws.Range("K45").Value = 1776.9231
$ws.Range("L45").Value = 1868.7142
$ws.Range("M45").Value = -1399.9231
$ws.Range("N45").Value = -2622.7142
$ws.Range("H61").Value = 1224.7709
$ws.Range("I61").Value = 983.13513
$ws.Range("J61").Value = 2037.5454
$ws.Range("K61").Value = 983.13513
$ws.Range("L61").Value = 2037.5454
$ws.Range("M61").Value = -771.13513
$ws.Range("N61").Value = -2461.5454
$ws.Range("H95").Value = 37016.8
$ws.Range("J95").Value = 37016.8
$ws.Range("L95").Value = 37016.8
$ws.Range("N95").Value = -42508.8
$ws.Range("H101").Value = 43453.2
$ws.Range("J101").Value = 43453.2
$ws.Range("L101").Value = 43453.2
$ws.Range("N101").Value = -49943.2
$ws.Range("H105").Value = 49244
$ws.Range("J105").Value = 49244
$ws.Range("L105").Value = 49244
$ws.Range("N105").Value = -56232
$ws.Range("H116").Value = 2205.4119
$ws.Range("I116").Value = 2466.5715
$ws.Range("J116").Value = 986.6667
$ws.Range("K116").Value = 2466.5715
$ws.Range("L116").Value = 986.6667
$ws.Range("M116").Value = -172.5715
$ws.Range("N116").Value = -5574.6667
$ws.Range("H122").Value = 1920.1875
$ws.Range("I122").Value = 1975.8334
$ws.Range("J122").Value = 1753.25
$ws.Range("K122").Value = 5927.5002
$ws.Range("L122").Value = 5259.75
$ws.Range("M122").Value = -3477.5002
$ws.Range("N122").Value = -10159.75
$ws.Range("H123").Value = 45950.332
$ws.Range("J123").Value = 45950.332
$ws.Range("L123").Value = 45950.332
$ws.Range("N123").Value = -55750.332
$ws.Range("H136").Value = 1224.7709
$ws.Range("I136").Value = 983.13513
$ws.Range("J136").Value = 2037.5454
$ws.Range("K136").Value = 2949.40539
$ws.Range("L136").Value = 6112.6362
$ws.Range("M136").Value = -399.4053899999999
$ws.Range("N136").Value = -11212.6362

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2205.4119
$ws.Range("I3").Value = 2466.5715
$ws.Range("J3").Value = 986.6667
$ws.Range("K3").Value = 2466.5715
$ws.Range("L3").Value = 986.6667
$ws.Range("M3").Value = -2352.5715
$ws.Range("N3").Value = -1214.6667
$ws.Range("H122").Value = 40673.6
$ws.Range("J122").Value = 40673.6
$ws.Range("L122").Value = 40673.6
$ws.Range("N122").Value = -50473.6

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3113.23
$ws.Range("I31").Value = 944.36365
$ws.Range("J31").Value = 3724.9614
$ws.Range("K31").Value = 944.36365
$ws.Range("L31").Value = 3724.9614
$ws.Range("M31").Value = -649.36365
$ws.Range("N31").Value = -4314.9614
$ws.Range("H34").Value = 3113.23
$ws.Range("I34").Value = 944.36365
$ws.Range("J34").Value = 3724.9614
$ws.Range("K34").Value = 944.36365
$ws.Range("L34").Value = 3724.9614
$ws.Range("M34").Value = -742.36365
$ws.Range("N34").Value = -4128.9614
$ws.Range("H43").Value = 40273
$ws.Range("J43").Value = 40273
$ws.Range("L43").Value = 40273
$ws.Range("N43").Value = -40641
$ws.Range("H101").Value = 40273
$ws.Range("J101").Value = 40273
$ws.Range("L101").Value = 40273
$ws.Range("N101").Value = -46763
$ws.Range("H131").Value = 41896
$ws.Range("J131").Value = 41896
$ws.Range("L131").Value = 41896
$ws.Range("N131").Value = -51976

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 3648.8
$ws.Range("I113").Value = 6184.778
$ws.Range("K113").Value = 18554.334
$ws.Range("M113").Value = -16384.334
$ws.Range("H133").Value = 7670
$ws.Range("I133").Value = 8005
$ws.Range("K133").Value = 24015
$ws.Range("M133").Value = -18955
$ws.Range("H134").Value = 91003624
$ws.Range("I134").Value = 91003624
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 273010872
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -273005802
$ws.Range("N134").ClearContents()
$ws.Range("H136").Value = 100001930
$ws.Range("I136").Value = 100001930
$ws.Range("K136").Value = 300005790
$ws.Range("M136").Value = -300000690
$ws.Range("H137").Value = 38471296
$ws.Range("J137").Value = 66681060
$ws.Range("L137").Value = 200043180
$ws.Range("N137").Value = -200053380
$ws.Range("H139").Value = 11896.154
$ws.Range("I139").Value = 13025
$ws.Range("J139").Value = 8133.3335
$ws.Range("K139").Value = 39075
$ws.Range("L139").Value = 24400.0005
$ws.Range("M139").Value = -33935
$ws.Range("N139").Value = -34680.00049999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 4221.091
$ws.Range("J92").Value = 4214.6665
$ws.Range("L92").Value = 4214.6665
$ws.Range("N92").Value = -7958.6665
$ws.Range("H97").Value = 6817.5
$ws.Range("I97").Value = 3726.25
$ws.Range("K97").Value = 3726.25
$ws.Range("M97").Value = -3230.25
$ws.Range("H104").Value = 47399.5
$ws.Range("J104").Value = 47399.5
$ws.Range("L104").Value = 47399.5
$ws.Range("N104").Value = -54387.5
$ws.Range("H113").Value = 1529.6364
$ws.Range("I113").Value = 1640
$ws.Range("K113").Value = 1640
$ws.Range("M113").Value = 530
$ws.Range("H124").Value = 41771.668
$ws.Range("J124").Value = 41771.668
$ws.Range("L124").Value = 41771.668
$ws.Range("N124").Value = -51591.668
$ws.Range("H126").Value = 10485.538
$ws.Range("I126").Value = 11392
$ws.Range("J126").Value = 5500
$ws.Range("K126").Value = 34176
$ws.Range("L126").Value = 16500
$ws.Range("M126").Value = -31706
$ws.Range("N126").Value = -21440

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H131").Value = 33250
$ws.Range("J131").Value = 33250
$ws.Range("L131").Value = 33250
$ws.Range("N131").Value = -43330
$ws.Range("H132").Value = 2523.4
$ws.Range("I132").Value = 1852.6531
$ws.Range("K132").Value = 5557.9593
$ws.Range("M132").Value = -3027.9593
$ws.Range("H136").Value = 1194.711
$ws.Range("I136").Value = 893.44446
$ws.Range("J136").Value = 2399.7778
$ws.Range("K136").Value = 2680.33338
$ws.Range("L136").Value = 7199.3334
$ws.Range("M136").Value = -130.33338
$ws.Range("N136").Value = -12299.3334

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 498.75
$ws.Range("I100").Value = 333
$ws.Range("K100").Value = 666
$ws.Range("M100").Value = -125
$ws.Range("H105").Value = 45538.332
$ws.Range("J105").Value = 45538.332
$ws.Range("L105").Value = 45538.332
$ws.Range("N105").Value = -52526.332
$ws.Range("H113").Value = 534.6875
$ws.Range("I113").Value = 550.15
$ws.Range("J113").Value = 508.91666
$ws.Range("K113").Value = 1650.45
$ws.Range("L113").Value = 1526.74998
$ws.Range("M113").Value = 519.5500000000002
$ws.Range("N113").Value = -5866.749980000001
$ws.Range("H122").Value = 7143507
$ws.Range("I122").Value = 9524343
$ws.Range("K122").Value = 28573029
$ws.Range("M122").Value = -28570579
$ws.Range("H132").Value = 1467.3208
$ws.Range("I132").Value = 1121.7556
$ws.Range("K132").Value = 3365.2668
$ws.Range("M132").Value = -835.2667999999999
